$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value in C2
$ws.Range("C2").Value = 12345

# Add new row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "H2JKV"
$ws.Range("C3").Value = 122333243

# Add new row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "H2JKV"
$ws.Range("C4").Value = 678568
